$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. "580.98").
# Excel auto-converts such text to a real number on assignment, but the
# source data keeps these as text (inline strings), matching the "Price"
# column format used throughout the sheet (e.g. "69.053.18", "3.472.70").
# Force text interpretation by pre-formatting as Text, then restore the
# original (default/"Normal") style so no formatting is left behind.
$textCells = @("D5", "D6", "D10", "D12", "D13", "D14", "D16", "D19", "D21", "D23", "D24", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D41", "D44", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin values (prices, volumes, names/links where ranking changed)
$ws.Range("D2").Value = '69.109.66'
$ws.Range("D3").Value = '3.475.19'
$ws.Range("E3").Value = '  -3.76%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '580.98'
$ws.Range("E5").Value = '  -0.81%  '
$ws.Range("D6").Value = '181.68'
$ws.Range("E6").Value = '  -4.86%  '
$ws.Range("D7").Value = '3.463.14'
$ws.Range("E7").Value = '  -4.03%  '
$ws.Range("E8").Value = '  -4.17%  '
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").Value = '0.195'
$ws.Range("E10").Value = '  +6.28%  '
$ws.Range("E11").Value = '  -4.16%  '
$ws.Range("D12").Value = '53.39'
$ws.Range("E12").Value = '  -4.97%  '
$ws.Range("D13").Value = '0.0000301'
$ws.Range("E13").Value = '  -4.13%  '
$ws.Range("D14").Value = '9.34'
$ws.Range("E14").Value = '  -4.24%  '
$ws.Range("D15").Value = '4.025.54'
$ws.Range("E15").Value = '  -3.90%  '
$ws.Range("D16").Value = '19.09'
$ws.Range("E16").Value = '  -4.63%  '
$ws.Range("D17").Value = '69.055.06'
$ws.Range("E17").Value = '  -1.92%  '
$ws.Range("D18").Value = '3.455.04'
$ws.Range("E18").Value = '  -4.20%  '
$ws.Range("D19").Value = '12.18'
$ws.Range("E19").Value = '  -4.39%  '
$ws.Range("E20").Value = '  -1.69%  '
$ws.Range("D21").Value = '537.07'
$ws.Range("E21").Value = '  +9.00%  '
$ws.Range("E22").Value = '  -4.80%  '
$ws.Range("D23").Value = '18.50'
$ws.Range("E23").Value = '  -7.62%  '
$ws.Range("D24").Value = '4.48'
$ws.Range("E24").Value = '  +2.19%  '
$ws.Range("E25").Value = '  -2.10%  '
$ws.Range("D26").Value = '94.67'
$ws.Range("E26").Value = '  -2.75%  '
$ws.Range("D27").Value = '11.10'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("E28").Value = '  -1.95%  '
$ws.Range("D29").Value = '9.02'
$ws.Range("E29").Value = '  -5.09%  '
$ws.Range("D30").Value = '31.43'
$ws.Range("E30").Value = '  -3.31%  '
$ws.Range("D31").Value = '7.16'
$ws.Range("E31").Value = '  -6.25%  '
$ws.Range("D32").Value = '12.40'
$ws.Range("E32").Value = '  +1.02%  '
$ws.Range("D33").Value = '63.52'
$ws.Range("E33").Value = '  -4.43%  '
$ws.Range("E34").Value = '  -6.22%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").Value = '3.10'
$ws.Range("E35").Value = '  +6.65%  '
$ws.Range("D36").Value = '0.404'
$ws.Range("E36").Value = '  +0.85%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").Value = '521.98'
$ws.Range("E37").Value = '  -9.96%  '
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").Value = '37.68'
$ws.Range("E39").Value = '  -3.58%  '
$ws.Range("D40").Value = '0.0₃0753'
$ws.Range("E40").Value = '  -8.10%  '
$ws.Range("D41").Value = '3.36'
$ws.Range("E41").Value = '  -3.12%  '
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '3.323.15'
$ws.Range("E43").Value = '  +2.96%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '3.07'
$ws.Range("E44").Value = '  -6.82%  '
$ws.Range("E45").Value = '  +2.73%  '
$ws.Range("E46").Value = '  -5.06%  '
$ws.Range("E48").Value = '  -4.07%  '
$ws.Range("D49").Value = '8.86'
$ws.Range("E49").Value = '  -8.37%  '
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").Value = '135.79'
$ws.Range("E51").Value = '  -2.31%  '

# Restore default styling on the text-forced cells (keeps them as text
# without leaving a lingering explicit number format applied).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
